# Add a new "FSAE_Achilles" worksheet as a copy of the existing
# "Sedan_HambaLG" template sheet, tweak its self-referential values,
# and make it the active tab.

$wb = $excel.ActiveWorkbook

# The existing template sheet.
$ws1 = $wb.Worksheets.Item(1)

# Duplicate it immediately after itself - this clones data, styles,
# column widths, frozen panes, and tab color in one shot.
$ws1.Copy($null, $ws1)

# The freshly copied sheet is now the 2nd tab; rename it.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "FSAE_Achilles"

# H3 on the template mirrors the sheet's own instance name - update it
# for the new sheet. H6 (rWheelCutout) gets a new value for this
# vehicle variant.
$ws2.Range("H3").Value = "FSAE_Achilles"
$ws2.Range("H6").Value = 0.25

# Make the new sheet the active/selected tab.
$ws2.Activate()
